# Insert a new data row at row 545 (shifts the existing rows 545:641 down to
# 546:642) and populate it with the new "Piña" / "Caramelo" price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(545).Insert()

$ws.Cells.Item(545, 1).Value = 10
$ws.Cells.Item(545, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(545, 3).Value = "La Araucanía"
$ws.Cells.Item(545, 4).Value = 44951
$ws.Cells.Item(545, 5).Value = 9
$ws.Cells.Item(545, 6).Value = "Fruta"
$ws.Cells.Item(545, 7).Value = 100108
$ws.Cells.Item(545, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(545, 9).Value = 100108005
$ws.Cells.Item(545, 10).Value = "Piña"
$ws.Cells.Item(545, 11).Value = "Caramelo"
$ws.Cells.Item(545, 12).Value = "Primera"
$ws.Cells.Item(545, 13).Value = 100
$ws.Cells.Item(545, 14).Value = 22000
$ws.Cells.Item(545, 15).Value = 22000
$ws.Cells.Item(545, 16).Value = 22000
$ws.Cells.Item(545, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item(545, 18).Value = "Ecuador"
$ws.Cells.Item(545, 19).Value = 1833
$ws.Cells.Item(545, 20).Value = 12

$ws.Cells.Item(545, 4).NumberFormat = $ws.Cells.Item(546, 4).NumberFormat
